# Added post-processed files for 80/20 clr cross validation
# Add 'hole_id' index column (A1:A32) to the 'train' worksheet so cross
# validation can be performed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
    "hole_id",
    "BRG_05_14",
    "BRG_16_02",
    "BRG_16_07",
    "BRG_05_12",
    "BRG_16_05",
    "ECO_09_02",
    "BRG_16_08",
    "BRG_16_04A",
    "BRG_01_05",
    "BRG_05_09",
    "BRG_05_02",
    "ECO_09_01",
    "BRG_13_02",
    "BRG_16_01",
    "BRG_16_09",
    "BRG_05_15",
    "BRG_05_13",
    "BRG_01_07",
    "BRG_05_05",
    "BRG_01_08",
    "ECO_09_05",
    "BRG_01_02",
    "BRG_05_11",
    "BRG_01_01",
    "BRG_05_10",
    "BRG_01_06",
    "BRG_05_03",
    "BRG_05_01",
    "BRG_01_09",
    "BRG_05_04",
    "ECO_09_04"
)

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}

# Cell A1 is brand new (row 1 previously only had headers in B1:M1), so it
# needs the same bold/centered/bordered header style as B1:M1. Copy the
# formatting (not the value, which we already set above) from B1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
